$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values in F and G columns for rows 348-449
$ws.Range("F348").Value = 232927
$ws.Range("F355").Value = 222155
$ws.Range("F360").Value = 749799
$ws.Range("F362").Value = 228997
$ws.Range("G362").Value = 3185
$ws.Range("F363").Value = 188724
$ws.Range("G363").Value = 2746
$ws.Range("F364").Value = 168742
$ws.Range("G364").Value = 2489
$ws.Range("F365").Value = 184631
$ws.Range("G365").Value = 2397
$ws.Range("F366").Value = 339457
$ws.Range("G366").Value = 2843
$ws.Range("F367").Value = 767039
$ws.Range("G367").Value = 3920
$ws.Range("F368").Value = 346241
$ws.Range("F370").Value = 180934
$ws.Range("G370").Value = 2044
$ws.Range("F371").Value = 160191
$ws.Range("G371").Value = 1966
$ws.Range("F372").Value = 178485
$ws.Range("G372").Value = 1854
$ws.Range("F373").Value = 350190
$ws.Range("G373").Value = 2378
$ws.Range("F374").Value = 773668
$ws.Range("G374").Value = 3422
$ws.Range("F375").Value = 351253
$ws.Range("G375").Value = 1851
$ws.Range("F376").Value = 222308
$ws.Range("G376").Value = 2228
$ws.Range("F378").Value = 157350
$ws.Range("G378").Value = 1548
$ws.Range("F379").Value = 179765
$ws.Range("G379").Value = 1613
$ws.Range("F380").Value = 344798
$ws.Range("G380").Value = 2025
$ws.Range("F381").Value = 746746
$ws.Range("G381").Value = 2692
$ws.Range("F382").Value = 356941
$ws.Range("F383").Value = 221255
$ws.Range("G383").Value = 1767
$ws.Range("F384").Value = 171968
$ws.Range("G384").Value = 1515
$ws.Range("F385").Value = 150846
$ws.Range("G385").Value = 1406
$ws.Range("F386").Value = 182894
$ws.Range("G386").Value = 1362
$ws.Range("F387").Value = 351595
$ws.Range("F388").Value = 730941
$ws.Range("G388").Value = 2206
$ws.Range("F389").Value = 353699
$ws.Range("F390").Value = 220014
$ws.Range("G390").Value = 1479
$ws.Range("F391").Value = 177595
$ws.Range("F392").Value = 221568
$ws.Range("G392").Value = 1214
$ws.Range("F397").Value = 107795
$ws.Range("F401").Value = 272418
$ws.Range("G401").Value = 936
$ws.Range("F409").Value = 708499
$ws.Range("G409").Value = 1007
$ws.Range("F415").Value = 307891
$ws.Range("F416").Value = 671641
$ws.Range("F418").Value = 202140
$ws.Range("F421").Value = 152971
$ws.Range("G421").Value = 532
$ws.Range("F423").Value = 439319
$ws.Range("F427").Value = 90399
$ws.Range("F428").Value = 102314
$ws.Range("F429").Value = 178348
$ws.Range("F432").Value = 123195
$ws.Range("F439").Value = 89013
$ws.Range("F440").Value = 73419
$ws.Range("F443").Value = 106520
$ws.Range("F444").Value = 103166
$ws.Range("F446").Value = 86221
$ws.Range("F447").Value = 67013
$ws.Range("F448").Value = 61080
$ws.Range("G448").Value = 136
$ws.Range("F449").Value = 58826
$ws.Range("G449").Value = 158

# Append new rows 450-452 with data through 2021-05-30
$ws.Range("A450").Value = 44344
$ws.Range("B450").Value = 389588
$ws.Range("C450").Value = 5589
$ws.Range("D450").Value = 148
$ws.Range("E450").Value = 12335
$ws.Range("F450").Value = 83116
$ws.Range("G450").Value = 161

$ws.Range("A451").Value = 44345
$ws.Range("B451").Value = 389690
$ws.Range("C451").Value = 2995
$ws.Range("D451").Value = 102
$ws.Range("E451").Value = 12339
$ws.Range("F451").Value = 76419
$ws.Range("G451").Value = 105

$ws.Range("A452").Value = 44346
$ws.Range("B452").Value = 389721
$ws.Range("C452").Value = 1309
$ws.Range("D452").Value = 31
$ws.Range("E452").Value = 12343
$ws.Range("F452").Value = 65231
$ws.Range("G452").Value = 111

